# Apply the commit's edits:
#  1. On the "Metadata" sheet, update URL, Version, Date, and Publisher values
#     to reflect the move from IBM/Alvearie to LinuxForHealth.
#  2. On the "Elements" sheet, clear the "Constraint(s)" cell for the
#     top-level "Extension" element row (it no longer carries the
#     ele-1/ext-1 constraint text -- that now only applies to
#     Extension.extension).

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/short-term-care-duration"
$wsMeta.Range("B3").Value = "8.0.0"
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

$wsElements = $wb.Worksheets.Item("Elements")
$wsElements.Range("AI2").Value = ""

# The "Extension.url" row's example value (column Q) shares the same
# underlying canonical URL text as the Metadata sheet's URL cell, so it
# must be kept in sync with the new URL as well.
$wsElements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/short-term-care-duration"

$wb.Save()
